# Update DividendHistory sheet: insert the latest dividend entry as a new
# second row (right under the header), pushing the existing history rows
# down by one. This mirrors a routine "new ex-dividend date" data refresh.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DividendHistory")

# Insert a new blank row above the current row 2 (the most-recent-date row),
# shifting all existing data rows down by one.
$ws.Rows.Item(2).Insert()

# Match the existing column formatting (values are stored as plain text,
# not as real dates/numbers), so force Text format before assigning.
$newRow = $ws.Range("A2:C2")
$newRow.NumberFormat = "@"

# Populate the newly inserted row with the latest dividend record.
$ws.Range("A2").Value = "03/10/2025"
$ws.Range("B2").Value = "03/10/2025"
$ws.Range("C2").Value = "0.010"
